$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.459.72'
$ws.Range('E2').Value = '  +0.07%  '

$ws.Range('D3').Value = '2.599.45'
$ws.Range('E3').Value = '  +6.94%  '

$ws.Range('D4').Formula = "'0.999"
$ws.Range('E4').Value = '  +0.09%  '

$ws.Range('D5').Formula = "'308.34"
$ws.Range('E5').Value = '  +4.03%  '

$ws.Range('D6').Formula = "'100.37"
$ws.Range('E6').Value = '  +2.98%  '

$ws.Range('D7').Formula = "'0.605"
$ws.Range('E7').Value = '  +5.86%  '

$ws.Range('E8').Value = '  +0.08%  '

$ws.Range('D9').Formula = "'0.580"
$ws.Range('E9').Value = '  +13.25%  '

$ws.Range('D10').Formula = "'39.35"
$ws.Range('E10').Value = '  +11.35%  '

$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').Formula = "'0.0847"
$ws.Range('E11').Value = '  +7.34%  '

$ws.Range('B12').Value = 'OKB'
$ws.Range('C12').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D12').Formula = "'54.43"
$ws.Range('E12').Value = '  +1.39%  '

$ws.Range('E13').Value = '  +14.06%  '

$ws.Range('D14').Value = '2.994.17'
$ws.Range('E14').Value = '  +6.86%  '

$ws.Range('E15').Value = '  +1.29%  '

$ws.Range('D16').Value = '2.602.79'
$ws.Range('E16').Value = '  +6.24%  '

$ws.Range('D17').Formula = "'0.926"
$ws.Range('E17').Value = '  +8.90%  '

$ws.Range('D18').Formula = "'15.05"
$ws.Range('E18').Value = '  +6.36%  '

$ws.Range('D19').Value = '46.574.11'
$ws.Range('E19').Value = '  +0.68%  '

$ws.Range('E20').Value = '  +7.09%  '

$ws.Range('D21').Formula = "'13.05"
$ws.Range('E21').Value = '  +2.40%  '

$ws.Range('D22').Formula = "'6.78"
$ws.Range('E22').Value = '  +8.62%  '

$ws.Range('D23').Formula = "'278.25"
$ws.Range('E23').Value = '  +13.22%  '

$ws.Range('D24').Formula = "'71.85"
$ws.Range('E24').Value = '  +6.23%  '

$ws.Range('D25').Formula = "'3.06"
$ws.Range('E25').Value = '  +9.01%  '

$ws.Range('D26').Formula = "'2.18"
$ws.Range('E26').Value = '  +11.18%  '

$ws.Range('D27').Formula = "'29.33"
$ws.Range('E27').Value = '  +37.01%  '

$ws.Range('E28').Value = '  +0.20%  '

$ws.Range('D29').Formula = "'4.04"
$ws.Range('E29').Value = '  +0.49%  '

$ws.Range('E30').Value = '  +8.87%  '

$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').Formula = "'39.13"
$ws.Range('E31').Value = '  -0.73%  '

$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').Formula = "'2.26"
$ws.Range('E32').Value = '  +1.36%  '

$ws.Range('D33').Formula = "'6.34"
$ws.Range('E33').Value = '  +13.88%  '

$ws.Range('D34').Formula = "'3.60"
$ws.Range('E34').Value = '  -6.44%  '

$ws.Range('E35').Value = '  +3.36%  '

$ws.Range('E36').Value = '  +8.73%  '

$ws.Range('D37').Formula = "'2.21"
$ws.Range('E37').Value = '  +9.42%  '

$ws.Range('D38').Formula = "'151.53"
$ws.Range('E38').Value = '  +2.20%  '

$ws.Range('E39').Value = '  +8.48%  '

$ws.Range('E40').Value = '  +5.82%  '

$ws.Range('D41').Formula = "'23.26"
$ws.Range('E41').Value = '  +41.22%  '

$ws.Range('E42').Value = '  +6.40%  '

$ws.Range('D43').Formula = "'0.0334"
$ws.Range('E43').Value = '  +9.92%  '

$ws.Range('D44').Formula = "'3.66"
$ws.Range('E44').Value = '  +11.69%  '

$ws.Range('D45').Formula = "'4.10"
$ws.Range('E45').Value = '  +4.14%  '

$ws.Range('D46').Value = '2.144.56'
$ws.Range('E46').Value = '  +8.16%  '

$ws.Range('E47').Value = '  -0.05%  '

$ws.Range('D48').Formula = "'93.10"
$ws.Range('E48').Value = '  +0.08%  '

$ws.Range('D49').Formula = "'9.49"
$ws.Range('E49').Value = '  +9.91%  '

$ws.Range('E50').Value = '  -2.84%  '

$ws.Range('D51').Formula = "'109.38"
$ws.Range('E51').Value = '  +7.97%  '
